$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 316
$ws.Range("I2").Value = 96.666664
$ws.Range("J2").Value = 645
$ws.Range("K2").Value = 96.666664
$ws.Range("L2").Value = 645
$ws.Range("M2").Value = 16.333336
$ws.Range("N2").Value = -871
$ws.Range("H113").Value = 12394.25
$ws.Range("I113").Value = 21847.5
$ws.Range("J113").Value = 2941
$ws.Range("K113").Value = 21847.5
$ws.Range("L113").Value = 2941
$ws.Range("M113").Value = -18593.5
$ws.Range("N113").Value = -9449
$ws.Range("H138").Value = 3588.303
$ws.Range("I138").Value = 2441.5715
$ws.Range("K138").Value = 7324.7145
$ws.Range("M138").Value = -2184.7145

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3450.5757
$ws.Range("I2").Value = 2468.739
$ws.Range("K2").Value = 2468.739
$ws.Range("M2").Value = -2355.739
$ws.Range("H3").Value = 325
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H45").Value = 1715.2106
$ws.Range("I45").Value = 1708
$ws.Range("J45").Value = 1735.4
$ws.Range("K45").Value = 1708
$ws.Range("L45").Value = 1735.4
$ws.Range("M45").Value = -1331
$ws.Range("N45").Value = -2489.4
$ws.Range("H105").Value = 47998
$ws.Range("J105").Value = 47998
$ws.Range("L105").Value = 47998
$ws.Range("N105").Value = -54986
$ws.Range("H110").Value = 2117410
$ws.Range("I110").Value = 2677624
$ws.Range("K110").Value = 2677624
$ws.Range("M110").Value = -2675579
$ws.Range("H116").Value = 3450.5757
$ws.Range("I116").Value = 2468.739
$ws.Range("K116").Value = 2468.739
$ws.Range("M116").Value = -174.739
$ws.Range("H122").Value = 2186117.8
$ws.Range("I122").Value = 4148512.2
$ws.Range("K122").Value = 12445536.6
$ws.Range("M122").Value = -12443086.6
$ws.Range("H132").Value = 2636093.5
$ws.Range("I132").Value = 1803.4615
$ws.Range("J132").Value = 8343722
$ws.Range("K132").Value = 5410.3845
$ws.Range("L132").Value = 25031166
$ws.Range("M132").Value = -2880.3845
$ws.Range("N132").Value = -25036226

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3450.5757
$ws.Range("I3").Value = 2468.739
$ws.Range("K3").Value = 2468.739
$ws.Range("M3").Value = -2354.739
$ws.Range("H86").Value = 58826628
$ws.Range("I86").Value = 3286.2856
$ws.Range("K86").Value = 3286.2856
$ws.Range("M86").Value = -2163.2856
$ws.Range("H89").Value = 58826628
$ws.Range("I89").Value = 3286.2856
$ws.Range("K89").Value = 16431.428
$ws.Range("M89").Value = -10815.428
$ws.Range("H107").Value = 979.7727
$ws.Range("I107").Value = 923.6875
$ws.Range("K107").Value = 923.6875
$ws.Range("M107").Value = 996.3125

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1374
$ws.Range("I2").Value = 1288.5
$ws.Range("K2").Value = 1288.5
$ws.Range("M2").Value = -1175.5
$ws.Range("H16").Value = 1019.6
$ws.Range("I16").Value = 1019.6
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1019.6
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -732.6
$ws.Range("N16").ClearContents()
$ws.Range("H32").Value = 2816.2856
$ws.Range("I32").Value = 2816.2856
$ws.Range("K32").Value = 2816.2856
$ws.Range("M32").Value = -2500.2856
$ws.Range("H33").Value = 2053.9
$ws.Range("I33").Value = 1557.75
$ws.Range("K33").Value = 1557.75
$ws.Range("M33").Value = -1178.75
$ws.Range("H113").Value = 1019.6
$ws.Range("I113").Value = 1019.6
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1019.6
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1150.4
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 40002290
$ws.Range("I132").Value = 2409.739
$ws.Range("J132").Value = 500000900
$ws.Range("K132").Value = 7229.217000000001
$ws.Range("L132").Value = 1500002700
$ws.Range("M132").Value = -4699.217000000001
$ws.Range("N132").Value = -1500007760
$ws.Range("H134").Value = 33339658
$ws.Range("I134").Value = 2290.913
$ws.Range("J134").Value = 142876720
$ws.Range("K134").Value = 6872.739
$ws.Range("L134").Value = 428630160
$ws.Range("M134").Value = -4337.739
$ws.Range("N134").Value = -428635230

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 310.64285
$ws.Range("I23").Value = 200.5
$ws.Range("K23").Value = 601.5
$ws.Range("M23").Value = -366.5
$ws.Range("H109").Value = 5557372
$ws.Range("I109").Value = 2179.8
$ws.Range("K109").Value = 6539.400000000001
$ws.Range("M109").Value = -5499.400000000001
$ws.Range("H131").Value = 1395.63
$ws.Range("I131").Value = 754.5714
$ws.Range("K131").Value = 2263.7142
$ws.Range("M131").Value = 2776.2858
$ws.Range("H132").Value = 1587
$ws.Range("I132").Value = 1404.5555
$ws.Range("J132").Value = 2056.1428
$ws.Range("K132").Value = 12640.9995
$ws.Range("L132").Value = 18505.2852
$ws.Range("M132").Value = -10110.9995
$ws.Range("N132").Value = -23565.2852

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 5030
$ws.Range("J49").Value = 5030
$ws.Range("L49").Value = 5030
$ws.Range("N49").Value = -5398
$ws.Range("H80").Value = 4429.6665
$ws.Range("J80").Value = 5671.2856
$ws.Range("L80").Value = 5671.2856
$ws.Range("N80").Value = -7667.2856
$ws.Range("H83").Value = 4429.6665
$ws.Range("J83").Value = 5671.2856
$ws.Range("L83").Value = 28356.428
$ws.Range("N83").Value = -38340.428
$ws.Range("H122").Value = 4243592.5
$ws.Range("I122").Value = 5657523.5
$ws.Range("J122").Value = 1799.5
$ws.Range("K122").Value = 16972570.5
$ws.Range("L122").Value = 5398.5
$ws.Range("M122").Value = -16970120.5
$ws.Range("N122").Value = -10298.5
$ws.Range("H126").Value = 9765975
$ws.Range("I126").Value = 5996026.5
$ws.Range("K126").Value = 17988079.5
$ws.Range("M126").Value = -17985609.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 1679.3334
$ws.Range("J13").Value = 1499
$ws.Range("L13").Value = 1499
$ws.Range("N13").Value = -1779
$ws.Range("H22").Value = 83334776
$ws.Range("I22").Value = 1625
$ws.Range("J22").Value = 250001070
$ws.Range("K22").Value = 1625
$ws.Range("L22").Value = 250001070
$ws.Range("M22").Value = -1330
$ws.Range("N22").Value = -250001660
$ws.Range("H27").Value = 83334776
$ws.Range("I27").Value = 1625
$ws.Range("J27").Value = 250001070
$ws.Range("K27").Value = 1625
$ws.Range("L27").Value = 250001070
$ws.Range("M27").Value = -1518
$ws.Range("N27").Value = -250001284
$ws.Range("H46").Value = 3360.647
$ws.Range("J46").Value = 3788.2856
$ws.Range("L46").Value = 3788.2856
$ws.Range("N46").Value = -4164.2856
$ws.Range("H74").Value = 47487.5
$ws.Range("I74").Value = 34975
$ws.Range("K74").Value = 34975
$ws.Range("M74").Value = -33977
$ws.Range("H77").Value = 47487.5
$ws.Range("I77").Value = 34975
$ws.Range("K77").Value = 104925
$ws.Range("M77").Value = -99933
$ws.Range("H82").Value = 1819.8422
$ws.Range("I82").Value = 2159.5715
$ws.Range("J82").Value = 868.6
$ws.Range("K82").Value = 2159.5715
$ws.Range("L82").Value = 868.6
$ws.Range("M82").Value = -1798.5715
$ws.Range("N82").Value = -1590.6
$ws.Range("H85").Value = 1819.8422
$ws.Range("I85").Value = 2159.5715
$ws.Range("J85").Value = 868.6
$ws.Range("K85").Value = 2159.5715
$ws.Range("L85").Value = 868.6
$ws.Range("M85").Value = -911.5715
$ws.Range("N85").Value = -3364.6
$ws.Range("H100").Value = 2535.25
$ws.Range("I100").Value = 2106.875
$ws.Range("K100").Value = 2106.875
$ws.Range("M100").Value = -1565.875
$ws.Range("H132").Value = 2233292.2
$ws.Range("I132").Value = 5799.7
$ws.Range("J132").Value = 5017658
$ws.Range("K132").Value = 17399.1
$ws.Range("L132").Value = 15052974
$ws.Range("M132").Value = -14869.1
$ws.Range("N132").Value = -15058034

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 22980.924
$ws.Range("I54").Value = 15915.272
$ws.Range("K54").Value = 15915.272
$ws.Range("M54").Value = -15395.272
$ws.Range("H81").Value = 525
$ws.Range("I81").Value = 525
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1050
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 11
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 525
$ws.Range("I84").Value = 525
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 5250
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 54
$ws.Range("N84").ClearContents()
$ws.Range("H132").Value = 6573.5
$ws.Range("I132").Value = 1348.2632
$ws.Range("J132").Value = 39666.668
$ws.Range("K132").Value = 4044.7896
$ws.Range("L132").Value = 119000.004
$ws.Range("M132").Value = -1514.7896
$ws.Range("N132").Value = -124060.004
